# JAN 1 Presenti Sheet
# Fill in the "Present"/"Absent" values for column D (Jan-02) on the Jan-2024 sheet
# for the first four students, and move the active cell selection to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jan-2024")

$ws.Range("D2").Value = "Present"
$ws.Range("D3").Value = "Present"
$ws.Range("D4").Value = "Absent"
$ws.Range("D5").Value = "Present"

$ws.Activate()
$ws.Range("D5").Select()
